$d = $word.ActiveDocument

# wdBrightGreen (4) serializes to OOXML w:highlight w:val="green"
$wdBrightGreen = 4

function Set-HighlightGreen($searchRange) {
    $searchRange.Font.HighlightColorIndex = $wdBrightGreen
}

function Find-FromStart($searchText, $wholeWord) {
    # Search the whole document for a (expected-unique) run of text and
    # return the matched Range, or $null if not found.
    $rng = $d.Content
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Text = $searchText
    $find.Forward = $true
    $find.Wrap = 0
    $find.MatchCase = $true
    $find.MatchWholeWord = $wholeWord
    $ok = $find.Execute()
    if ($ok) {
        return $rng
    }
    return $null
}

function Find-FromPosition($searchText, $startPos, $wholeWord) {
    # Search starting at a given character offset (used to disambiguate
    # separator text like " / " or ", " that repeats elsewhere in the doc).
    $rng = $d.Range($startPos, $d.Content.End)
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Text = $searchText
    $find.Forward = $true
    $find.Wrap = 0
    $find.MatchCase = $true
    $find.MatchWholeWord = $wholeWord
    $ok = $find.Execute()
    if ($ok) {
        return $rng
    }
    return $null
}

# --- "get, recv" row: only the "recv" run changes -----------------------
$r = Find-FromStart "recv" $true
if ($r -ne $null) { Set-HighlightGreen $r }

# --- "ascii / binary" row: "ascii", " / ", and "binary" all change ------
$r = Find-FromStart "ascii" $true
if ($r -ne $null) {
    Set-HighlightGreen $r
    $pos = $r.End

    $r2 = Find-FromPosition " / " $pos $false
    if ($r2 -ne $null) {
        Set-HighlightGreen $r2
        $pos = $r2.End
    }

    $r3 = Find-FromPosition "binary" $pos $true
    if ($r3 -ne $null) { Set-HighlightGreen $r3 }
}

# --- "open, close" row: only the "close" run changes ---------------------
$r = Find-FromStart "close" $true
if ($r -ne $null) { Set-HighlightGreen $r }

# --- "quit, bye" row: only the "bye" run changes --------------------------
$r = Find-FromStart "bye" $true
if ($r -ne $null) { Set-HighlightGreen $r }

# --- "help, ?" row: "help", ", ", and "?" all change ----------------------
$r = Find-FromStart "help" $true
if ($r -ne $null) {
    Set-HighlightGreen $r
    $pos = $r.End

    $r2 = Find-FromPosition ", " $pos $false
    if ($r2 -ne $null) {
        Set-HighlightGreen $r2
        $pos = $r2.End
    }

    $r3 = Find-FromPosition "?" $pos $false
    if ($r3 -ne $null) { Set-HighlightGreen $r3 }
}

Write-Output "Done applying highlight changes"
